$d = $word.ActiveDocument

# --- First paragraph ---
# Originally: "**ID__AFFARS_5312_topic_7__ID** " (the id-placeholder run,
# followed by a trailing single-space run). Target: a single run containing
# "**ID__AFFARS_SUBPART_5312_3__ID**" (no trailing space run).
$oldId = "**ID__AFFARS_5312_topic_7__ID**"
$newId = "**ID__AFFARS_SUBPART_5312_3__ID**"

$para1Range = $d.Paragraphs(1).Range
$pStart = $para1Range.Start
$pEnd = $para1Range.End

# Drop everything after the id text but before the paragraph mark (the
# trailing space run) first, then rewrite the id run's text in place so its
# rPr is left untouched.
$d.Range($pStart + $oldId.Length, $pEnd - 1).Delete()
$d.Range($pStart, $pStart + $oldId.Length).Text = $newId

# Add a paragraph border (5pt space on every side) and widen the left indent.
$para1 = $d.Paragraphs(1)
$para1.Format.Borders.DistanceFromTop = 5
$para1.Format.Borders.DistanceFromLeft = 5
$para1.Format.Borders.DistanceFromBottom = 5
$para1.Format.Borders.DistanceFromRight = 5
$para1.Format.LeftIndent = 11.25
